$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H37").Value = 1702.5
$ws.Range("I37").Value = 1055
$ws.Range("J37").Value = 2350
$ws.Range("K37").Value = 3165
$ws.Range("L37").Value = 7050
$ws.Range("M37").Value = -3039
$ws.Range("N37").Value = -7302

$ws.Range("H43").Value = 2286.0833
$ws.Range("I43").Value = 3654.8
$ws.Range("J43").Value = 1308.4286
$ws.Range("K43").Value = 3654.8
$ws.Range("L43").Value = 1308.4286
$ws.Range("M43").Value = -3585.8
$ws.Range("N43").Value = -1446.4286

$ws.Range("H50").Value = 337
$ws.Range("J50").Value = 337
$ws.Range("L50").Value = 1011
$ws.Range("N50").Value = -1961

$ws.Range("H59").Value = 931.875
$ws.Range("I59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("M59").ClearContents()

$ws.Range("H94").Value = 6987.5
$ws.Range("I94").Value = 984.7143
$ws.Range("J94").Value = 49007
$ws.Range("K94").Value = 984.7143
$ws.Range("L94").Value = 49007
$ws.Range("M94").Value = -533.7143
$ws.Range("N94").Value = -49909

$ws.Range("H98").Value = 4541.6665
$ws.Range("I98").Value = 4524.0386
$ws.Range("K98").Value = 4524.0386
$ws.Range("M98").Value = -3026.0386

$ws.Range("H99").Value = 805.6
$ws.Range("I99").Value = 845.2857
$ws.Range("J99").Value = 250
$ws.Range("K99").Value = 2535.8571
$ws.Range("L99").Value = 750
$ws.Range("M99").Value = -1037.8571
$ws.Range("N99").Value = -3746

$ws.Range("H100").Value = 1413.2
$ws.Range("I100").Value = 1299.75
$ws.Range("J100").Value = 1867
$ws.Range("K100").Value = 1299.75
$ws.Range("L100").Value = 1867
$ws.Range("M100").Value = -758.75
$ws.Range("N100").Value = -2949

$ws.Range("H103").Value = 516.1667
$ws.Range("I103").Value = 449.5
$ws.Range("J103").Value = 649.5
$ws.Range("K103").Value = 1348.5
$ws.Range("L103").Value = 1948.5
$ws.Range("M103").Value = -762.5
$ws.Range("N103").Value = -3120.5

$ws.Range("H106").Value = 5323.222
$ws.Range("J106").Value = 17899.334
$ws.Range("L106").Value = 17899.334
$ws.Range("N106").Value = -19161.334

$ws.Range("H113").Value = 11335.111
$ws.Range("I113").Value = 11335.111
$ws.Range("K113").Value = 11335.111
$ws.Range("M113").Value = -8081.111000000001

$ws.Range("H116").Value = 3708063
$ws.Range("I116").Value = 4682796
$ws.Range("J116").Value = 4079
$ws.Range("K116").Value = 4682796
$ws.Range("L116").Value = 4079
$ws.Range("M116").Value = -4679354
$ws.Range("N116").Value = -10963

$ws.Range("H122").Value = 4541.6665
$ws.Range("I122").Value = 4524.0386
$ws.Range("K122").Value = 13572.1158
$ws.Range("M122").Value = -11122.1158

$ws.Range("H132").Value = 10525.424
$ws.Range("I132").Value = 8117.6
$ws.Range("J132").Value = 34603.668
$ws.Range("K132").Value = 24352.8
$ws.Range("L132").Value = 103811.004
$ws.Range("M132").Value = -21822.8
$ws.Range("N132").Value = -108871.004

$ws.Range("H137").Value = 5517.68
$ws.Range("I137").Value = 2523.8027
$ws.Range("J137").Value = 12847.518
$ws.Range("K137").Value = 7571.408100000001
$ws.Range("L137").Value = 38542.554
$ws.Range("M137").Value = -5021.408100000001
$ws.Range("N137").Value = -43642.554

$ws.Range("H138").Value = 6304.12
$ws.Range("I138").Value = 17000
$ws.Range("J138").Value = 5374.0435
$ws.Range("K138").Value = 51000
$ws.Range("L138").Value = 16122.1305
$ws.Range("M138").Value = -45860
$ws.Range("N138").Value = -26402.1305

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3399.2876
$ws.Range("I32").Value = 2142.1384
$ws.Range("J32").Value = 13613.625
$ws.Range("K32").Value = 2142.1384
$ws.Range("L32").Value = 13613.625
$ws.Range("M32").Value = -1855.1384
$ws.Range("N32").Value = -14187.625

$ws.Range("H102").Value = 2535.2856
$ws.Range("I102").Value = 2535.2856
$ws.Range("K102").Value = 2535.2856
$ws.Range("M102").Value = -913.2856000000002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1568.9445
$ws.Range("I94").Value = 902.5625
$ws.Range("J94").Value = 6900
$ws.Range("K94").Value = 902.5625
$ws.Range("L94").Value = 6900
$ws.Range("M94").Value = -451.5625
$ws.Range("N94").Value = -7802

$ws.Range("H99").Value = 16747.191
$ws.Range("I99").Value = 16696.643
$ws.Range("J99").Value = 16848.285
$ws.Range("K99").Value = 16696.643
$ws.Range("L99").Value = 16848.285
$ws.Range("M99").Value = -15198.643
$ws.Range("N99").Value = -19844.285

$ws.Range("H105").Value = 1513.3334
$ws.Range("I105").Value = 1395
$ws.Range("J105").Value = 1750
$ws.Range("K105").Value = 1395
$ws.Range("L105").Value = 1750
$ws.Range("M105").Value = 352
$ws.Range("N105").Value = -5244

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4591.5884
$ws.Range("J16").Value = 13036.75
$ws.Range("L16").Value = 13036.75
$ws.Range("N16").Value = -13610.75

$ws.Range("H22").Value = 2700
$ws.Range("I22").Value = 2700
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 2700
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -2350
$ws.Range("N22").ClearContents()

$ws.Range("H31").Value = 13327.782
$ws.Range("I31").Value = 5071.1763
$ws.Range("J31").Value = 17021.525
$ws.Range("K31").Value = 5071.1763
$ws.Range("L31").Value = 17021.525
$ws.Range("M31").Value = -4776.1763
$ws.Range("N31").Value = -17611.525

$ws.Range("H34").Value = 13327.782
$ws.Range("I34").Value = 5071.1763
$ws.Range("J34").Value = 17021.525
$ws.Range("K34").Value = 5071.1763
$ws.Range("L34").Value = 17021.525
$ws.Range("M34").Value = -4869.1763
$ws.Range("N34").Value = -17425.525

$ws.Range("H113").Value = 4591.5884
$ws.Range("J113").Value = 13036.75
$ws.Range("L113").Value = 13036.75
$ws.Range("N113").Value = -17376.75

$ws.Range("H122").Value = 3236.0435
$ws.Range("I122").Value = 1181.9375
$ws.Range("J122").Value = 7931.143
$ws.Range("K122").Value = 3545.8125
$ws.Range("L122").Value = 23793.429
$ws.Range("M122").Value = -1095.8125
$ws.Range("N122").Value = -28693.429

$ws.Range("H129").Value = 87500
$ws.Range("J129").Value = 87500
$ws.Range("L129").Value = 87500
$ws.Range("N129").Value = -97500

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 4171.027
$ws.Range("I68").Value = 3124.5
$ws.Range("J68").Value = 4297.879
$ws.Range("K68").Value = 9373.5
$ws.Range("L68").Value = 12893.637
$ws.Range("M68").Value = -8562.5
$ws.Range("N68").Value = -14515.637

$ws.Range("H71").Value = 4171.027
$ws.Range("I71").Value = 3124.5
$ws.Range("J71").Value = 4297.879
$ws.Range("K71").Value = 28120.5
$ws.Range("L71").Value = 38680.911
$ws.Range("M71").Value = -24064.5
$ws.Range("N71").Value = -46792.911

$ws.Range("H122").Value = 12663246
$ws.Range("I122").Value = 23360752
$ws.Range("J122").Value = 3154351.8
$ws.Range("K122").Value = 210246768
$ws.Range("L122").Value = 28389166.2
$ws.Range("M122").Value = -210244318
$ws.Range("N122").Value = -28394066.2

$ws.Range("H131").Value = 1451.32
$ws.Range("J131").Value = 1486.8617
$ws.Range("L131").Value = 4460.5851
$ws.Range("N131").Value = -14540.5851

$ws.Range("H132").Value = 2756602.8
$ws.Range("I132").Value = 1700
$ws.Range("K132").Value = 15300
$ws.Range("M132").Value = -12770

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 24080.834
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 24080.834
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 24080.834
$ws.Range("M52").ClearContents()
$ws.Range("N52").Value = -24598.834

$ws.Range("H97").Value = 4463.143
$ws.Range("I97").Value = 1500
$ws.Range("J97").Value = 5271.273
$ws.Range("K97").Value = 1500
$ws.Range("L97").Value = 5271.273
$ws.Range("M97").Value = -1004
$ws.Range("N97").Value = -6263.273

$ws.Range("H122").Value = 13664.363
$ws.Range("I122").Value = 8923.223
$ws.Range("J122").Value = 34999.5
$ws.Range("K122").Value = 26769.669
$ws.Range("L122").Value = 104998.5
$ws.Range("M122").Value = -24319.669
$ws.Range("N122").Value = -109898.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1045.75
$ws.Range("I46").Value = 871.3333
$ws.Range("J46").Value = 2266.6667
$ws.Range("K46").Value = 871.3333
$ws.Range("L46").Value = 2266.6667
$ws.Range("M46").Value = -683.3333
$ws.Range("N46").Value = -2642.6667

$ws.Range("H61").Value = 2256.4285
$ws.Range("I61").Value = 2321.84
$ws.Range("J61").Value = 1711.3334
$ws.Range("K61").Value = 2321.84
$ws.Range("L61").Value = 1711.3334
$ws.Range("M61").Value = -2119.84
$ws.Range("N61").Value = -2115.3334

$ws.Range("H100").Value = 8521.071
$ws.Range("I100").Value = 6471.143
$ws.Range("J100").Value = 10571
$ws.Range("K100").Value = 6471.143
$ws.Range("L100").Value = 10571
$ws.Range("M100").Value = -5930.143
$ws.Range("N100").Value = -11653

$ws.Range("H113").Value = 2256.4285
$ws.Range("I113").Value = 2321.84
$ws.Range("J113").Value = 1711.3334
$ws.Range("K113").Value = 2321.84
$ws.Range("L113").Value = 1711.3334
$ws.Range("M113").Value = -151.8400000000001
$ws.Range("N113").Value = -6051.3334

$ws.Range("H139").Value = 94900
$ws.Range("J139").Value = 94900
$ws.Range("L139").Value = 94900
$ws.Range("N139").Value = -105180
